$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ISO canton-code values for column K, row 1 (header) through row 163,
# in the same order as the worksheet rows.
$isoValues = @(
    "iso",
    "CH",
    "VD",
    "VS",
    "GE",
    "BE",
    "FR",
    "SO",
    "NE",
    "JU",
    "BS",
    "BL",
    "AG",
    "ZH",
    "GL",
    "SH",
    "AR",
    "AI",
    "SG",
    "GR",
    "TG",
    "LU",
    "UR",
    "SZ",
    "OW",
    "NW",
    "ZG",
    "TI",
    "CH",
    "VD",
    "VS",
    "GE",
    "BE",
    "FR",
    "SO",
    "NE",
    "JU",
    "BS",
    "BL",
    "AG",
    "ZH",
    "GL",
    "SH",
    "AR",
    "AI",
    "SG",
    "GR",
    "TG",
    "LU",
    "UR",
    "SZ",
    "OW",
    "NW",
    "ZG",
    "TI",
    "CH",
    "VD",
    "VS",
    "GE",
    "BE",
    "FR",
    "SO",
    "NE",
    "JU",
    "BS",
    "BL",
    "AG",
    "ZH",
    "GL",
    "SH",
    "AR",
    "AI",
    "SG",
    "GR",
    "TG",
    "LU",
    "UR",
    "SZ",
    "OW",
    "NW",
    "ZG",
    "TI",
    "CH",
    "VD",
    "VS",
    "GE",
    "BE",
    "FR",
    "SO",
    "NE",
    "JU",
    "BS",
    "BL",
    "AG",
    "ZH",
    "GL",
    "SH",
    "AR",
    "AI",
    "SG",
    "GR",
    "TG",
    "LU",
    "UR",
    "SZ",
    "OW",
    "NW",
    "ZG",
    "TI",
    "CH",
    "VD",
    "VS",
    "GE",
    "BE",
    "FR",
    "SO",
    "NE",
    "JU",
    "BS",
    "BL",
    "AG",
    "ZH",
    "GL",
    "SH",
    "AR",
    "AI",
    "SG",
    "GR",
    "TG",
    "LU",
    "UR",
    "SZ",
    "OW",
    "NW",
    "ZG",
    "TI",
    "CH",
    "VD",
    "VS",
    "GE",
    "BE",
    "FR",
    "SO",
    "NE",
    "JU",
    "BS",
    "BL",
    "AG",
    "ZH",
    "GL",
    "SH",
    "AR",
    "AI",
    "SG",
    "GR",
    "TG",
    "LU",
    "UR",
    "SZ",
    "OW",
    "NW",
    "ZG",
    "TI"
)

for ($i = 0; $i -lt $isoValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 11).Value = $isoValues[$i]
}

# K1 (header "iso") should carry the same formatting as the other header
# cells in row 1 (e.g. J1).
$ws.Cells.Item(1, 11).NumberFormat = $ws.Cells.Item(1, 10).NumberFormat

# Reflect the workbook's on-screen selection state at the time of the edit.
[void]$ws.Range("K137:K163").Select()
